# "After Hours" tracklist refresh: weeknd5 -> weeknd6
# - replaces the 7-track placeholder data with the real 14-track "After Hours" listing
# - drops the old "Composer" column header, folds it into "Title/Composer"
# - renames the weeknd5 defined names to weeknd6 and grows their range to row 15
# - Sheet2's report rows recompute automatically off the Sheet1 formulas

$wb = $excel.ActiveWorkbook

$tracks = @(
    @(1,  "Alone Again",                  "Adam Feeney / Carlo `"Illangelo`" Montagnese / Jason Quenneville / Abel Tesfaye", 0.17500000000000002),
    @(2,  "Too Late",                     "Eric Frederic / Carlo `"Illangelo`" Montagnese / Jason Quenneville / Abel Tesfaye", 0.16597222222222222),
    @(3,  "Hardest to Love",              "Oscar Holter / Max Martin / Abel Tesfaye", 0.14652777777777778),
    @(4,  "Scared to Live",               "Ahmad Balshe / Oscar Holter / Elton John / Daniel Lopatin / Max Martin / Bernie Taupin / Abel Tesfaye", 0.13263888888888889),
    @(5,  "Snowchild",                    "Ahmad Balshe / Carlo `"Illangelo`" Montagnese / Jason Quenneville / Abel Tesfaye", 0.17152777777777775),
    @(6,  "Escape from LA",               "Metro Boomin / Mike McTaggart / Carlo `"Illangelo`" Montagnese / Abel Tesfaye", 0.24652777777777779),
    @(7,  "Heartless",                    "Metro Boomin / Carlo `"Illangelo`" Montagnese / Andre Proctor / Abel Tesfaye", 0.1388888888888889),
    @(8,  "Faith",                        "Ahmad Balshe / Metro Boomin / Carlo `"Illangelo`" Montagnese / Abel Tesfaye", 0.19652777777777777),
    @(9,  "Blinding Lights",              "Ahmad Balshe / Oscar Holter / Max Martin / Jason Quenneville / Abel Tesfaye", 0.14166666666666666),
    @(10, "In Your Eyes",                 "Ahmad Balshe / Oscar Holter / Max Martin / Abel Tesfaye", 0.16458333333333333),
    @(11, "Save Your Tears",              "Ahmad Balshe / Oscar Holter / Max Martin / Jason Quenneville / Abel Tesfaye", 0.14930555555555555),
    @(12, "Repeat After Me (Interlude)",  "Daniel Lopatin / Kevin Parker / Abel Tesfaye", 0.13541666666666666),
    @(13, "After Hours",                  "Ahmad Balshe / Carlo `"Illangelo`" Montagnese / Jason Quenneville / Abel Tesfaye / Mario Winans", 0.25138888888888888),
    @(14, "Until I Bleed Out",            "Metro Boomin / Daniel Lopatin / Notinbed / Mejdi Rhars / Abel Tesfaye", 0.13194444444444445)
)

foreach ($sheetName in @("Sheet1", "Sheet3")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Header row: the separate "Composer" header column goes away, title
    # header becomes "Title/Composer"; performer/time headers lose their
    # trailing-space padding.
    $ws.Range("B1").Value = "Title/Composer"
    $ws.Range("C1").ClearContents()
    $ws.Range("D1").Value = "Performer"
    $ws.Range("E1").Value = "Time"

    foreach ($t in $tracks) {
        $row = 1 + $t[0]
        $ws.Cells.Item($row, 1).Value = $t[0]
        $ws.Cells.Item($row, 2).Value = $t[1]
        $ws.Cells.Item($row, 3).Value = $t[2]
        $ws.Cells.Item($row, 4).Value = "The Weeknd"
        $ws.Cells.Item($row, 5).Value = $t[3]
    }
}

# Defined names: weeknd5 -> weeknd6, range grows from 13 to 15 rows
$n1 = $wb.Names.Item("Sheet1!weeknd5")
$n1.Name = "weeknd6"
$n1b = $wb.Names.Item("Sheet1!weeknd6")
$n1b.RefersTo = "=Sheet1!`$A`$1:`$E`$15"

$n3 = $wb.Names.Item("Sheet3!weeknd5")
$n3.Name = "weeknd6"
$n3b = $wb.Names.Item("Sheet3!weeknd6")
$n3b.RefersTo = "=Sheet3!`$A`$1:`$E`$15"

# Sheet2's report view now spans down to row 18 (14 tracks + 4 header/blank rows)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A3:K18").Select()
